# qpsk 2400 rx wip
# Update input parameters on "LoopFilter LPF" and "Branch LPF" sheets, then
# leave "Branch LPF" as the active/selected sheet (tabSelected + activeTab),
# with the cell selections left where the author last clicked.

$wb = $excel.ActiveWorkbook

# --- LoopFilter LPF (sheet1) ---
$ws1 = $wb.Worksheets.Item("LoopFilter LPF")
$ws1.Range("B2").Value = 50       # cutoff freq: 100 -> 50 Hz
$ws1.Range("B3").Value = 14400    # sample rate: 28800 -> 14400 Hz
$ws1.Range("B16").Value = 8       # gain: 1 -> 8

# --- Branch LPF (sheet3) ---
$ws3 = $wb.Worksheets.Item("Branch LPF")
$ws3.Range("B2").Value = 1200     # cutoff freq: 4800 -> 1200 Hz
$ws3.Range("B3").Value = 14400    # sample rate: 28800 -> 14400 Hz
$ws3.Range("B18").Value = 1       # gain: 2 -> 1

# Leave the selection on LoopFilter LPF parked at B17 (no longer the active
# tab), then activate Branch LPF and select B28 there, matching the saved
# view state in the workbook.
$ws1.Activate()
$ws1.Range("B17").Select()

$ws3.Activate()
$ws3.Range("B28").Select()
